$wb = $excel.ActiveWorkbook

# ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1339.1428
$ws.Range("I12").Value = 793.3333
$ws.Range("J12").Value = 1748.5
$ws.Range("K12").Value = 793.3333
$ws.Range("L12").Value = 1748.5
$ws.Range("M12").Value = -623.3333
$ws.Range("N12").Value = -2088.5

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1416.15
$ws.Range("I19").Value = 831.3333
$ws.Range("K19").Value = 831.3333
$ws.Range("M19").Value = -656.3333

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8666.333000000001
$ws.Range("J32").Value = 9999.5
$ws.Range("L32").Value = 9999.5
$ws.Range("N32").Value = -10651.5

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 471.92307
$ws.Range("I33").Value = 346.63635
$ws.Range("K33").Value = 346.63635
$ws.Range("M33").Value = -117.63635

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2781.2307
$ws.Range("I38").Value = 1442.6666
$ws.Range("J38").Value = 3928.5715
$ws.Range("K38").Value = 4327.9998
$ws.Range("L38").Value = 11785.7145
$ws.Range("M38").Value = -3955.9998
$ws.Range("N38").Value = -12529.7145

# ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 112.333336
$ws.Range("I39").Value = 110.35294
$ws.Range("K39").Value = 331.05882
$ws.Range("M39").Value = -35.05882000000003

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4036.4
$ws.Range("I40").Value = 4599.5
$ws.Range("J40").Value = 3895.625
$ws.Range("K40").Value = 4599.5
$ws.Range("L40").Value = 3895.625
$ws.Range("M40").Value = -4424.5
$ws.Range("N40").Value = -4245.625

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3737.3794
$ws.Range("I132").Value = 3513.7144
$ws.Range("K132").Value = 10541.1432
$ws.Range("M132").Value = -8011.143199999999

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 9004.354499999999
$ws.Range("I137").Value = 9985.925999999999
$ws.Range("J137").Value = 2378.75
$ws.Range("K137").Value = 29957.778
$ws.Range("L137").Value = 7136.25
$ws.Range("M137").Value = -27407.778
$ws.Range("N137").Value = -12236.25

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3655.2827
$ws.Range("I138").Value = 1461.2727
$ws.Range("K138").Value = 4383.8181
$ws.Range("M138").Value = 756.1818999999996

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6369.143
$ws.Range("I141").Value = 5651.636
$ws.Range("K141").Value = 16954.908
$ws.Range("M141").Value = -11774.908

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6987.2144
$ws.Range("I45").Value = 6771
$ws.Range("J45").Value = 7203.4287
$ws.Range("K45").Value = 6771
$ws.Range("L45").Value = 7203.4287
$ws.Range("M45").Value = -6394
$ws.Range("N45").Value = -7957.4287

# ARM row 62
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248

# ARM row 65
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240

# ARM row 94
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H94").Value = 188590000
$ws.Range("J94").Value = 220009170
$ws.Range("L94").Value = 220009170
$ws.Range("N94").Value = -220010972

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6200.4
$ws.Range("I86").Value = 5182.7144
$ws.Range("J86").Value = 8575
$ws.Range("K86").Value = 5182.7144
$ws.Range("L86").Value = 8575
$ws.Range("M86").Value = -4059.7144
$ws.Range("N86").Value = -10821

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 6200.4
$ws.Range("I89").Value = 5182.7144
$ws.Range("J89").Value = 8575
$ws.Range("K89").Value = 25913.572
$ws.Range("L89").Value = 42875
$ws.Range("M89").Value = -20297.572
$ws.Range("N89").Value = -54107

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 16757.46
$ws.Range("I7").Value = 30913.572
$ws.Range("K7").Value = 30913.572
$ws.Range("M7").Value = -30800.572

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 706.75
$ws.Range("I22").Value = 586.2857
$ws.Range("J22").Value = 1550
$ws.Range("K22").Value = 586.2857
$ws.Range("L22").Value = 1550
$ws.Range("M22").Value = -236.2857
$ws.Range("N22").Value = -2250

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4384.5356
$ws.Range("I134").Value = 2216.5789
$ws.Range("J134").Value = 8961.333000000001
$ws.Range("K134").Value = 6649.736699999999
$ws.Range("L134").Value = 26883.999
$ws.Range("M134").Value = -4114.736699999999
$ws.Range("N134").Value = -31953.999

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 338.25
$ws.Range("I23").Value = 21
$ws.Range("J23").Value = 383.57144
$ws.Range("K23").Value = 63
$ws.Range("L23").Value = 1150.71432
$ws.Range("M23").Value = 172
$ws.Range("N23").Value = -1620.71432

# CUL row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 35252.89
$ws.Range("I97").Value = 56378.184
$ws.Range("K97").Value = 169134.552
$ws.Range("M97").Value = -168638.552

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 44038
$ws.Range("J132").Value = 65624.875
$ws.Range("L132").Value = 590623.875
$ws.Range("N132").Value = -595683.875

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 9021.647000000001
$ws.Range("I133").Value = 3708.3333
$ws.Range("J133").Value = 14999.125
$ws.Range("K133").Value = 11124.9999
$ws.Range("L133").Value = 44997.375
$ws.Range("M133").Value = -6064.999899999999
$ws.Range("N133").Value = -55117.375

# GSM row 92
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 15601.429
$ws.Range("J92").Value = 15601.429
$ws.Range("L92").Value = 15601.429
$ws.Range("N92").Value = -19345.429

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 726.0476
$ws.Range("J107").Value = 718.75
$ws.Range("L107").Value = 718.75
$ws.Range("N107").Value = -4558.75

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12095.529
$ws.Range("J122").Value = 13453.125
$ws.Range("L122").Value = 40359.375
$ws.Range("N122").Value = -45259.375

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 14588.4375
$ws.Range("I126").Value = 13584.6875
$ws.Range("K126").Value = 40754.0625
$ws.Range("M126").Value = -38284.0625

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3575.4285

# LTW row 13
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 4899
$ws.Range("I13").Value = 4899
$ws.Range("K13").Value = 4899
$ws.Range("M13").Value = -4759

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2133
$ws.Range("I100").Value = 2133
$ws.Range("K100").Value = 2133
$ws.Range("M100").Value = -1592

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3953.348
$ws.Range("I136").Value = 1461.75
$ws.Range("J136").Value = 7829.1665
$ws.Range("K136").Value = 4385.25
$ws.Range("L136").Value = 23487.4995
$ws.Range("M136").Value = -1835.25
$ws.Range("N136").Value = -28587.4995

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 28515
$ws.Range("I123").Value = 28515
$ws.Range("K123").Value = 28515
$ws.Range("M123").Value = -23615
